# Update to analysis 2 based on error report from evan
# The "outcomes" sheet (sheet2.xml) had a row (cohort_definition_id 745,
# "Inflammatory Bowel Disease") that was a duplicate/erroneous entry.
# Remove that row entirely; all subsequent rows shift up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("outcomes")

# Row 57 currently holds cohort_definition_id = 745 ("Inflammatory Bowel
# Disease"), which is an erroneous duplicate row. Delete the entire row so
# everything below shifts up by one row (ids/names/clean windows keep their
# original association, only their row position changes).
$ws.Rows.Item(57).Delete()
